# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.233.61"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.66%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.071.25"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.88%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "529.59"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +6.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.95"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +7.17%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +5.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.114"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +7.90%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +6.42%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.591.43"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.40"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +9.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000174"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +17.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.152.96"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.23"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +8.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.075.05"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.19"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +7.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.23"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +6.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "342.84"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +5.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.507"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +8.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.57"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +5.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0₃0979"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +9.85%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +4.73%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.03"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +9.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.50"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +10.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.87"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +7.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.24"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +6.48%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.84"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +9.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "158.07"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.00"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +7.76%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.23"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +13.51%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.102.98"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.86"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.99"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +12.67%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.83%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.49"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +6.40%  "
$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.05"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +5.71%  "
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.339.60"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.03"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.09"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +6.20%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.31"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +6.98%  "
